$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect so the cells below can be edited, then
# restore protection afterwards.
$ws.Unprotect()

# Update the confidential disclaimer date (2021-05-13 -> 2021-05-14)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) columns for rows 2-8
$ws.Range("D2").Value = 0.5028562690505626
$ws.Range("E2").Value = 0.01243411271793482

$ws.Range("D3").Value = 0.2419643238123078
$ws.Range("E3").Value = 0.01845184518451837

$ws.Range("D4").Value = 0.09487752130641366
$ws.Range("E4").Value = 0.01846153846153853

$ws.Range("D5").Value = 0.1037783225522941
$ws.Range("E5").Value = 0.01386046511627903

$ws.Range("D6").Value = 0.02990553373046482
$ws.Range("E6").Value = 0.02101491380980058

$ws.Range("D7").Value = 0.02661802954795707
$ws.Range("E7").Value = 0.02401659618606877

$ws.Range("E8").Value = 0.01517499728245597

# Restore sheet protection
$ws.Protect()
